# "Generate Report for Handoff": refresh the handoff report for the files
# that were just handed off in this run. For each of those rows, the
# Priority column (on both target-language sheets) picks up the "ht"
# (handoff-type) marker, and the "Latest Handoff Datetime" / "Latest HO
# Xliff Generate Date" timestamps move forward a few seconds to reflect the
# new report-generation time: zh-cn's handoff pass finished at 18:23:47,
# de-de's finished a few seconds later at 18:23:52, and the Overview sheet
# mirrors the de-de (last-finished) timestamp.

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 10, 11, 13, 14)

# Target-language sheets: mark these rows as handed-off ("ht") and bump the
# "Latest Handoff Datetime" (column H) to each sheet's new generation time.
$ws = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value = "ht"
    $ws.Cells.Item($r, 8).Value = "2016-09-05 18:23:47"
}

$ws = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value = "ht"
    $ws.Cells.Item($r, 8).Value = "2016-09-05 18:23:52"
}

# Overview sheet: bump "Latest HO Xliff Generate Date" (column G) for the
# same set of rows - it mirrors the de-de handoff timestamp.
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-09-05 18:23:52"
}
